# Update numeric values across three worksheets to match the new
# distractor-analysis output (see commit: "Added unit tests for utils
# and distractor analysis.").

$wb = $excel.ActiveWorkbook

# --- Sheet "gof": recompute Deviance / AIC / BIC for the two rows ---
$wsGof = $wb.Worksheets.Item("gof")
$wsGof.Range("D2").Value = 15637
$wsGof.Range("F2").Value = 15673
$wsGof.Range("G2").Value = 15775

$wsGof.Range("D3").Value = 15628
$wsGof.Range("F3").Value = 15694
$wsGof.Range("G3").Value = 15880

# --- Sheet "Estimates 0-1": xsi / se.xsi / std / Femp per item row ---
$wsEst = $wb.Worksheets.Item("Estimates 0-1")

$wsEst.Range("B2").Value = -0.089
$wsEst.Range("D2").Value = -0.084
$wsEst.Range("E2").Value = 1.408

$wsEst.Range("B3").Value = -0.082
$wsEst.Range("D3").Value = -0.077
$wsEst.Range("E3").Value = 1.591

$wsEst.Range("B4").Value = -0.107
$wsEst.Range("C4").Value = 0.085
$wsEst.Range("D4").Value = -0.101
$wsEst.Range("E4").Value = 1.585

$wsEst.Range("B5").Value = 0.05
$wsEst.Range("C5").Value = 0.1
$wsEst.Range("D5").Value = 0.047
$wsEst.Range("E5").Value = 0.25

$wsEst.Range("B6").Value = -0.022
$wsEst.Range("D6").Value = -0.021
$wsEst.Range("E6").Value = 0.05

$wsEst.Range("B7").Value = -0.042
$wsEst.Range("D7").Value = -0.04
$wsEst.Range("E7").Value = 0.195

$wsEst.Range("B8").Value = 0.123
$wsEst.Range("C8").Value = 0.094
$wsEst.Range("D8").Value = 0.116
$wsEst.Range("E8").Value = 1.712

$wsEst.Range("B9").Value = -0.168
$wsEst.Range("C9").Value = 0.094
$wsEst.Range("D9").Value = -0.159
$wsEst.Range("E9").Value = 3.194

$wsEst.Range("B10").Value = -0.044
$wsEst.Range("C10").Value = 0.093
$wsEst.Range("D10").Value = -0.042
$wsEst.Range("E10").Value = 0.224

$wsEst.Range("B11").Value = -0.182
$wsEst.Range("D11").Value = -0.172
$wsEst.Range("E11").Value = 3.749

$wsEst.Range("B12").Value = 0.125
$wsEst.Range("C12").Value = 0.097
$wsEst.Range("D12").Value = 0.118
$wsEst.Range("E12").Value = 1.661

$wsEst.Range("B13").Value = -0.173
$wsEst.Range("D13").Value = -0.163
$wsEst.Range("E13").Value = 2.934

$wsEst.Range("B14").Value = 0.181
$wsEst.Range("D14").Value = 0.171
$wsEst.Range("E14").Value = 3.788

$wsEst.Range("B15").Value = 0.201
$wsEst.Range("C15").Value = 0.094
$wsEst.Range("D15").Value = 0.19
$wsEst.Range("E15").Value = 4.572

$wsEst.Range("B16").Value = 0.323
$wsEst.Range("D16").Value = 0.305
$wsEst.Range("E16").Value = 10.863

$wsEst.Range("B17").Value = 0.094
$wsEst.Range("D17").Value = 0.089
$wsEst.Range("E17").Value = 0.069

# --- Sheet "Main effect 0-1": Unstandardized / Standardized estimates ---
$wsMain = $wb.Worksheets.Item("Main effect 0-1")

$wsMain.Range("B2").Value = -0.413
$wsMain.Range("C2").Value = -0.39

$wsMain.Range("B3").Value = -0.4
$wsMain.Range("C3").Value = -0.377
